# Update the "Förändrad" (Changed) date column (C) for rows 2-15
# from 2023-10-25 (serial 45224) to 2023-11-03 (serial 45233).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
